# Updates cryptos list values (price/volume/coin rearrangements) per source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "53.642.75"
$ws.Range("E2").Value = "  -4.37%  "

# Row 3
$ws.Range("D3").Value = "2.224.91"
$ws.Range("E3").Value = "  -5.97%  "

# Row 4
$ws.Range("E4").Value = "  -0.10%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "484.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.24%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "125.52"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.80%  "

# Row 7
$ws.Range("E7").Value = "  -0.06%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.517"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.04%  "

# Row 9
$ws.Range("D9").Value = "2.230.59"
$ws.Range("E9").Value = "  -5.70%  "

# Row 10
$ws.Range("E10").Value = "  -6.55%  "

# Row 11
$ws.Range("E11").Value = "  -1.38%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.65"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.78%  "

# Row 13
$ws.Range("E13").Value = "  -2.65%  "

# Row 14
$ws.Range("D14").Value = "2.620.95"
$ws.Range("E14").Value = "  -5.99%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.98"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.91%  "

# Row 16
$ws.Range("D16").Value = "53.551.03"
$ws.Range("E16").Value = "  -4.51%  "

# Row 17
$ws.Range("E17").Value = "  -3.36%  "

# Row 18
$ws.Range("D18").Value = "2.232.90"
$ws.Range("E18").Value = "  -4.63%  "

# Row 19
$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.95"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.60%  "

# Row 20
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.55"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.42%  "

# Row 21
$ws.Range("E21").Value = "  -2.53%  "

# Row 22
$ws.Range("E22").Value = "  -2.48%  "

# Row 23
$ws.Range("E23").Value = "  -0.08%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.29"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.02%  "

# Row 25
$ws.Range("E25").Value = "  +0.01%  "

# Row 26
$ws.Range("E26").Value = "  -1.38%  "

# Row 27
$ws.Range("E27").Value = "  -3.16%  "

# Row 28
$ws.Range("E28").Value = "  -3.35%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "169.90"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.71%  "

# Row 30
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.57"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.68%  "

# Row 31
$ws.Range("B31").Value = "PEPE"
$ws.Range("C31").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D31").Value = "0.0₃0675"
$ws.Range("E31").Value = "  -4.99%  "

# Row 32
$ws.Range("E32").Value = "  -0.15%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.996"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.13%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.71"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.27%  "

# Row 35
$ws.Range("E35").Value = "  -3.31%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.40"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.67%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.15"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.84%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.830"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.04%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.57"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.57%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "35.84"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.46%  "

# Row 41
$ws.Range("E41").Value = "  -0.70%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.36"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.63%  "

# Row 43
$ws.Range("E43").Value = "  -2.15%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "122.52"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.62%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.61"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.96%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0876"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.97%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.533"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.87%  "

# Row 48
$ws.Range("B48").Value = "Bittensor"
$ws.Range("C48").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "228.87"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.30%  "

# Row 49
$ws.Range("B49").Value = "Hedera"
$ws.Range("C49").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0468"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.42%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0201"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.73%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "15.93"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.22%  "
